# Commit: "solved sort the jumbled numbers"
# Adds a new tracker row (row 24) for the LeetCode problem
# "Sort the Jumbled Numbers", including its hyperlink, and updates the
# current selection to the newly added cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Populate the new row's values first (column order C,D,E,F,G) so the
#        shared-string table grows in the same order as the source workbook. ---
$ws.Range("A24").Value = 2191
$ws.Range("B24").Value = "Medium"
$ws.Range("C24").Value = "Sort the Jumbled Numbers"
$ws.Range("D24").Value = "https://leetcode.com/problems/sort-the-jumbled-numbers/description/"
$ws.Range("E24").Value = "Array"
$ws.Range("F24").Value = "O(nd+logn)"
$ws.Range("G24").Value = "Encode the numbers according to the map, store (encoded_num, index) in a list. Exploit python's sort feature. Check the notes, this is good problem."
$ws.Range("H24").Value = 45497

# --- 2. Copy the formatting (fill/border/alignment banding) from the row
#        directly above so the new row matches the table's existing style. ---
$ws.Range("A23:H23").Copy()
$ws.Range("A24:H24").PasteSpecial(-4122)

# --- 3. Wire up the hyperlink for the new question's URL cell, then restore
#        that cell's banded/bordered style (Hyperlinks.Add resets it). ---
$ws.Hyperlinks.Add($ws.Range("D24"), "https://leetcode.com/problems/sort-the-jumbled-numbers/description/")
$ws.Range("D23").Copy()
$ws.Range("D24").PasteSpecial(-4122)

# --- 4. Match the author's final on-screen selection. ---
$ws.Range("G24").Select()
